# resolved issues of RW05 and RW06
# The "SmokeTest" flag (column D) on the ScenarioMapping sheet was
# incorrectly set on the RW05 scenarios (rows 99-110) and missing on the
# RW06 scenarios (rows 111-122). Flip them: RW05 -> "No", RW06 -> "Yes".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioMapping")

# RW05 scenarios (TC_098 .. TC_109) - turn SmokeTest off
$ws.Range("D99:D110").Value = "No"

# RW06 scenarios (TC_110 .. TC_121) - turn SmokeTest on
$ws.Range("D111:D122").Value = "Yes"

# Leave the view focused on the area that was just edited, matching the
# last cell touched (RW06's first row).
$ws.Range("C111").Select()
